# Update the Correspond Handoff Datetime / Correspond Handback DateTime
# values on the zh-cn and de-de report sheets (regenerating the report
# produced new timestamps for the d690e023... row on each sheet).

$wb = $excel.ActiveWorkbook

# zh-cn sheet: row 3 holds the d690e023-... entry
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D3").Value = "2016-01-11 05:04:42"
$wsZhCn.Range("G3").Value = "2016-01-11 05:05:47"

# de-de sheet: row 3 holds the d690e023-... entry
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D3").Value = "2016-01-11 05:04:58"
$wsDeDe.Range("G3").Value = "2016-01-11 05:06:14"
